$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 600
$ws.Range("J32").Value = 600
$ws.Range("L32").Value = 600
$ws.Range("N32").Value = -1252
$ws.Range("H64").Value = 3500
$ws.Range("I64").Value = 3300
$ws.Range("J64").Value = 3572.7273
$ws.Range("K64").Value = 3300
$ws.Range("L64").Value = 3572.7273
$ws.Range("M64").Value = -3052
$ws.Range("N64").Value = -4068.7273
$ws.Range("H67").Value = 3500
$ws.Range("I67").Value = 3300
$ws.Range("J67").Value = 3572.7273
$ws.Range("K67").Value = 3300
$ws.Range("L67").Value = 3572.7273
$ws.Range("M67").Value = -2442
$ws.Range("N67").Value = -5288.7273
$ws.Range("H135").Value = 1258.1578
$ws.Range("I135").Value = 966.94446
$ws.Range("J135").Value = 6500
$ws.Range("K135").Value = 8702.50014
$ws.Range("L135").Value = 58500
$ws.Range("M135").Value = -6167.50014
$ws.Range("N135").Value = -63570
$ws.Range("H138").Value = 7577061
$ws.Range("J138").Value = 2684.75
$ws.Range("L138").Value = 8054.25
$ws.Range("N138").Value = -18334.25
$ws.Range("H140").Value = 44666.668
$ws.Range("J140").Value = 44666.668
$ws.Range("L140").Value = 44666.668
$ws.Range("N140").Value = -55026.668

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1102.0625
$ws.Range("I2").Value = 620.125
$ws.Range("J2").Value = 1584
$ws.Range("K2").Value = 620.125
$ws.Range("L2").Value = 1584
$ws.Range("M2").Value = -507.125
$ws.Range("N2").Value = -1810
$ws.Range("H37").Value = 7867
$ws.Range("I37").Value = 3400.5
$ws.Range("J37").Value = 16800
$ws.Range("K37").Value = 3400.5
$ws.Range("L37").Value = 16800
$ws.Range("M37").Value = -3127.5
$ws.Range("N37").Value = -17346
$ws.Range("H44").Value = 22700
$ws.Range("J44").Value = 22700
$ws.Range("L44").Value = 22700
$ws.Range("N44").Value = -23676
$ws.Range("H55").Value = 24400
$ws.Range("J55").Value = 24400
$ws.Range("L55").Value = 24400
$ws.Range("N55").Value = -25030
$ws.Range("H63").Value = 1668100
$ws.Range("I63").Value = 2501325
$ws.Range("J63").Value = 1650
$ws.Range("K63").Value = 2501325
$ws.Range("L63").Value = 1650
$ws.Range("M63").Value = -2500639
$ws.Range("N63").Value = -3022
$ws.Range("H66").Value = 1668100
$ws.Range("I66").Value = 2501325
$ws.Range("J66").Value = 1650
$ws.Range("K66").Value = 12506625
$ws.Range("L66").Value = 8250
$ws.Range("M66").Value = -12503193
$ws.Range("N66").Value = -15114
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H116").Value = 1102.0625
$ws.Range("I116").Value = 620.125
$ws.Range("J116").Value = 1584
$ws.Range("K116").Value = 620.125
$ws.Range("L116").Value = 1584
$ws.Range("M116").Value = 1673.875
$ws.Range("N116").Value = -6172

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1102.0625
$ws.Range("I3").Value = 620.125
$ws.Range("J3").Value = 1584
$ws.Range("K3").Value = 620.125
$ws.Range("L3").Value = 1584
$ws.Range("M3").Value = -506.125
$ws.Range("N3").Value = -1812
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 5000
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 5000
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 5000
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -5348
$ws.Range("H31").Value = 4169251.2
$ws.Range("I31").Value = 1684.6333
$ws.Range("J31").Value = 16671951
$ws.Range("K31").Value = 1684.6333
$ws.Range("L31").Value = 16671951
$ws.Range("M31").Value = -1389.6333
$ws.Range("N31").Value = -16672541
$ws.Range("H34").Value = 4169251.2
$ws.Range("I34").Value = 1684.6333
$ws.Range("J34").Value = 16671951
$ws.Range("K34").Value = 1684.6333
$ws.Range("L34").Value = 16671951
$ws.Range("M34").Value = -1482.6333
$ws.Range("N34").Value = -16672355
$ws.Range("H41").Value = 6750
$ws.Range("I41").Value = 1500
$ws.Range("J41").Value = 12000
$ws.Range("K41").Value = 1500
$ws.Range("L41").Value = 12000
$ws.Range("M41").Value = -1072
$ws.Range("N41").Value = -12856
$ws.Range("H42").Value = 12996
$ws.Range("I42").Value = 10000
$ws.Range("J42").Value = 13745
$ws.Range("K42").Value = 10000
$ws.Range("L42").Value = 13745
$ws.Range("M42").Value = -9407
$ws.Range("N42").Value = -14931
$ws.Range("H50").Value = 12344.8
$ws.Range("J50").Value = 13160.889
$ws.Range("L50").Value = 13160.889
$ws.Range("N50").Value = -14410.889
$ws.Range("H51").Value = 24666.666
$ws.Range("J51").Value = 24666.666
$ws.Range("L51").Value = 24666.666
$ws.Range("N51").Value = -26138.666
$ws.Range("H59").Value = 44666.668
$ws.Range("J59").Value = 44666.668
$ws.Range("L59").Value = 44666.668
$ws.Range("N59").Value = -46956.668
$ws.Range("H60").Value = 16974.555
$ws.Range("J60").Value = 16974.555
$ws.Range("L60").Value = 16974.555
$ws.Range("N60").Value = -17996.555
$ws.Range("H61").Value = 24666.666
$ws.Range("J61").Value = 24666.666
$ws.Range("L61").Value = 24666.666
$ws.Range("N61").Value = -25362.666
$ws.Range("H74").Value = 29782.8
$ws.Range("J74").Value = 29782.8
$ws.Range("L74").Value = 29782.8
$ws.Range("N74").Value = -31530.8
$ws.Range("H77").Value = 29782.8
$ws.Range("J77").Value = 29782.8
$ws.Range("L77").Value = 89348.39999999999
$ws.Range("N77").Value = -98084.39999999999
$ws.Range("H132").Value = 2354.5898
$ws.Range("I132").Value = 1638.4814
$ws.Range("J132").Value = 3965.8333
$ws.Range("K132").Value = 4915.4442
$ws.Range("L132").Value = 11897.4999
$ws.Range("M132").Value = -2385.4442
$ws.Range("N132").Value = -16957.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1174.95
$ws.Range("J5").Value = 1549.0769
$ws.Range("L5").Value = 4647.2307
$ws.Range("N5").Value = -4871.2307
$ws.Range("H131").Value = 337592.16
$ws.Range("J131").Value = 529989
$ws.Range("L131").Value = 1589967
$ws.Range("N131").Value = -1600047
$ws.Range("H132").Value = 2173.7778
$ws.Range("I132").Value = 1018.9091
$ws.Range("J132").Value = 3988.5715
$ws.Range("K132").Value = 9170.1819
$ws.Range("L132").Value = 35897.1435
$ws.Range("M132").Value = -6640.1819
$ws.Range("N132").Value = -40957.1435
$ws.Range("H135").Value = 1174.95
$ws.Range("J135").Value = 1549.0769
$ws.Range("L135").Value = 13941.6921
$ws.Range("N135").Value = -19011.6921

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 5400
$ws.Range("I43").Value = 800
$ws.Range("K43").Value = 800
$ws.Range("M43").Value = -649
$ws.Range("H46").Value = 14900
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 14900
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 14900
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -15212
$ws.Range("H57").Value = 10900
$ws.Range("J57").Value = 19800
$ws.Range("L57").Value = 19800
$ws.Range("N57").Value = -21440
$ws.Range("H80").Value = 2366.3333
$ws.Range("I80").Value = 2300
$ws.Range("J80").Value = 2399.5
$ws.Range("K80").Value = 2300
$ws.Range("L80").Value = 2399.5
$ws.Range("M80").Value = -1302
$ws.Range("N80").Value = -4395.5
$ws.Range("H83").Value = 2366.3333
$ws.Range("I83").Value = 2300
$ws.Range("J83").Value = 2399.5
$ws.Range("K83").Value = 11500
$ws.Range("L83").Value = 11997.5
$ws.Range("M83").Value = -6508
$ws.Range("N83").Value = -21981.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2246.6667
$ws.Range("I22").Value = 3070
$ws.Range("K22").Value = 3070
$ws.Range("M22").Value = -2775
$ws.Range("H27").Value = 2246.6667
$ws.Range("I27").Value = 3070
$ws.Range("K27").Value = 3070
$ws.Range("M27").Value = -2963
$ws.Range("H46").Value = 1612.125
$ws.Range("J46").Value = 1612.125
$ws.Range("L46").Value = 1612.125
$ws.Range("N46").Value = -1988.125
$ws.Range("H47").Value = 17750
$ws.Range("J47").Value = 17750
$ws.Range("L47").Value = 17750
$ws.Range("N47").Value = -18730
$ws.Range("H52").Value = 17750
$ws.Range("J52").Value = 17750
$ws.Range("L52").Value = 17750
$ws.Range("N52").Value = -18216
